$d = $word.ActiveDocument

# 1. Replace the existing sentence text with the new, shorter sentence.
$d.Content.Find.Execute(
    "De header wordt maar een keer benoemd in het bericht. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "De header attributen worden als HTTP headers verstuurd.", 2) | Out-Null

# 2. Find the paragraph that now holds that sentence.
$p = $null
$pIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq "De header attributen worden als HTTP headers verstuurd.") {
        $p = $d.Paragraphs($i)
        $pIndex = $i
        break
    }
}

# 3. Insert a brand-new list paragraph right after it (Word copies the
#    numPr/spacing/rPr of the split point onto the new paragraph).
$p.Range.InsertParagraphAfter()
$newp = $d.Paragraphs($pIndex + 1)
$newp.Range.Text = "De header wordt maar een keer benoemd in de YAML-specificatie. "

# 4. Split "de YAML-specificatie" into its own run, and park a
#    zero-width "_GoBack" bookmark right after it (i.e. before the
#    closing ". "). A document can only have one bookmark with a given
#    name, so adding "_GoBack" here automatically relocates the one
#    that previously sat in the "De specificatie heeft een aantal
#    uitgangspunten:" paragraph.
$r = $newp.Range
$prefixLen = "De header wordt maar een keer benoemd in ".Length
$midLen = "de YAML-specificatie".Length
$splitStart = $r.Start + $prefixLen
$splitEnd = $splitStart + $midLen

# Throwaway bookmark forces a run break before "de YAML-specificatie"
# without touching any run-formatting (rPr stays identical across runs).
$tempRange = $d.Range($splitStart, $splitStart)
$d.Bookmarks.Add("ZzTempSplit", $tempRange) | Out-Null

$goBackRange = $d.Range($splitEnd, $splitEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

$d.Bookmarks("ZzTempSplit").Delete()

Write-Output "edit complete"
